# Update the "Price" column (D) for the refreshed coin symbol list.
# Values must remain text (matching the original inlineStr cells) so that
# formatting such as trailing zeros (e.g. "0.09320") and non-exponential
# small numbers (e.g. "0.00000000751") is preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "271.88"
    3  = "23.03"
    4  = "6.366"
    5  = "0.06297"
    6  = "3.656"
    7  = "6.757"
    8  = "1.402"
    9  = "0.8352"
    10 = "0.1623"
    11 = "0.08404"
    13 = "0.03122"
    14 = "0.09308"
    15 = "3.962"
    16 = "0.001707"
    17 = "0.04855"
    18 = "0.006226"
    19 = "0.005476"
    21 = "0.0001500"
    22 = "3.735"
    23 = "2.326"
    24 = "0.01386"
    25 = "0.3379"
    26 = "0.1218"
    27 = "0.0002682"
    40 = "0.04683"
    41 = "0.006890"
    43 = "0.003459"
    45 = "0.00006244"
    46 = "0.00000000750"
    47 = "0.7887"
    48 = "0.1151"
    49 = "0.00002100"
    50 = "0.01240"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    # Force the cell to be treated as text so numeric-looking strings are
    # not coerced into floating point numbers (which would lose trailing
    # zeros / switch to scientific notation).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    # Drop the temporary "Text" number format so no extra style is left
    # behind on the cell (matches original unstyled cells).
    $cell.Style = "Normal"
}

$wb.Save()
